# Rename AHB-Diff column headers to carry the format-version suffix
# (_old -> _FV2410, _new -> _FV2504), turn the header row + data range into
# a proper Excel Table, and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header cells (row 1) -----------------------------------
$oldCols  = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newCols  = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")
$baseNames = @(
    "Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID",
    "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung"
)

for ($i = 0; $i -lt $oldCols.Length; $i++) {
    $ws.Range("$($oldCols[$i])1").Value = "$($baseNames[$i])_FV2410"
}
for ($i = 0; $i -lt $newCols.Length; $i++) {
    $ws.Range("$($newCols[$i])1").Value = "$($baseNames[$i])_FV2504"
}
# Column K ("diff") keeps its text as-is.

# --- 2. Turn A1:U57 into an Excel Table -----------------------------------
# Stash a copy of the header row's current formatting on a scratch row far
# away, strip formatting from the real header row (so the engine doesn't
# synthesize a headerRowDxfId/dxf override when the table is created), add
# the table, then paste the stashed formatting back onto the header row and
# wipe the scratch row again. Net effect: cell styles end up unchanged.
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy() | Out-Null
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null

$headerRange.ClearFormats()

$dataRange = $ws.Range("A1:U57")
$tbl = $ws.ListObjects.Add(
    [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
    $dataRange,
    $null,
    [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratchRange.Copy() | Out-Null
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats) | Out-Null
$scratchRange.Clear() | Out-Null

# --- 3. Freeze the header row ----------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
